$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd28"
$ws.Range("C2").Value = "Cd80"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1791996666666667
$ws.Range("H2").Value = 0.537599
$ws.Range("I2").Value = 0.04251079199666429
$ws.Range("J2").Value = 0.04251079199666429
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 21.557693
$ws.Range("N2").Value = 64.673079
$ws.Range("O2").Value = 0.5505707555812251
$ws.Range("P2").Value = 0.5505707555812251
$ws.Range("Q2").Value = 3.863131399702334
$ws.Range("R2").Value = 34.768182597321
$ws.Range("S2").Value = 0.02340519886995976
$ws.Range("T2").Value = 0.02340519886995976

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd28"
$ws.Range("C3").Value = "Cd80"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1791996666666667
$ws.Range("H3").Value = 0.537599
$ws.Range("I3").Value = 0.04251079199666429
$ws.Range("J3").Value = 0.04251079199666429
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.981869333333333
$ws.Range("N3").Value = 11.945608
$ws.Range("O3").Value = 0.1016945926207894
$ws.Range("P3").Value = 0.1016945926207894
$ws.Range("Q3").Value = 0.7135496572435556
$ws.Range("R3").Value = 6.421946915192001
$ws.Range("S3").Value = 0.004323117674087888
$ws.Range("T3").Value = 0.004323117674087888

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd28"
$ws.Range("C4").Value = "Cd80"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1791996666666667
$ws.Range("H4").Value = 0.537599
$ws.Range("I4").Value = 0.04251079199666429
$ws.Range("J4").Value = 0.04251079199666429
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.00696933333334
$ws.Range("N4").Value = 36.02090800000001
$ws.Range("O4").Value = 0.3066509100994217
$ws.Range("P4").Value = 0.3066509100994217
$ws.Range("Q4").Value = 2.151644902210223
$ws.Range("R4").Value = 19.364804119892
$ws.Range("S4").Value = 0.01303597305482432
$ws.Range("T4").Value = 0.01303597305482432

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cd28"
$ws.Range("C5").Value = "Cd80"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1791996666666667
$ws.Range("H5").Value = 0.537599
$ws.Range("I5").Value = 0.04251079199666429
$ws.Range("J5").Value = 0.04251079199666429
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.608641
$ws.Range("N5").Value = 4.825923
$ws.Range("O5").Value = 0.04108374169856382
$ws.Range("P5").Value = 0.04108374169856382
$ws.Range("Q5").Value = 0.2882679309863333
$ws.Range("R5").Value = 2.594411378877
$ws.Range("S5").Value = 0.00174650239779233
$ws.Range("T5").Value = 0.00174650239779233

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cd28"
$ws.Range("C6").Value = "Cd80"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.036192666666667
$ws.Range("H6").Value = 12.108578
$ws.Range("I6").Value = 0.9574892080033358
$ws.Range("J6").Value = 0.9574892080033357
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 21.557693
$ws.Range("N6").Value = 64.673079
$ws.Range("O6").Value = 0.5505707555812251
$ws.Range("P6").Value = 0.5505707555812251
$ws.Range("Q6").Value = 87.01100239685134
$ws.Range("R6").Value = 783.099021571662
$ws.Range("S6").Value = 0.5271655567112654
$ws.Range("T6").Value = 0.5271655567112654

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd28"
$ws.Range("C7").Value = "Cd80"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.036192666666667
$ws.Range("H7").Value = 12.108578
$ws.Range("I7").Value = 0.9574892080033358
$ws.Range("J7").Value = 0.9574892080033357
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.981869333333333
$ws.Range("N7").Value = 11.945608
$ws.Range("O7").Value = 0.1016945926207894
$ws.Range("P7").Value = 0.1016945926207894
$ws.Range("Q7").Value = 16.07159180282489
$ws.Range("R7").Value = 144.644326225424
$ws.Range("S7").Value = 0.09737147494670148
$ws.Range("T7").Value = 0.09737147494670148

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd28"
$ws.Range("C8").Value = "Cd80"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.036192666666667
$ws.Range("H8").Value = 12.108578
$ws.Range("I8").Value = 0.9574892080033358
$ws.Range("J8").Value = 0.9574892080033357
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.00696933333334
$ws.Range("N8").Value = 36.02090800000001
$ws.Range("O8").Value = 0.3066509100994217
$ws.Range("P8").Value = 0.3066509100994217
$ws.Range("Q8").Value = 48.46244157209156
$ws.Range("R8").Value = 436.1619741488241
$ws.Range("S8").Value = 0.2936149370445975
$ws.Range("T8").Value = 0.2936149370445974

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd28"
$ws.Range("C9").Value = "Cd80"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.036192666666667
$ws.Range("H9").Value = 12.108578
$ws.Range("I9").Value = 0.9574892080033358
$ws.Range("J9").Value = 0.9574892080033357
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.608641
$ws.Range("N9").Value = 4.825923
$ws.Range("O9").Value = 0.04108374169856382
$ws.Range("P9").Value = 0.04108374169856382
$ws.Range("Q9").Value = 6.492785007499332
$ws.Range("R9").Value = 58.43506506749399
$ws.Range("S9").Value = 0.0393372393007715
$ws.Range("T9").Value = 0.03933723930077149
